# Update the cryptos price list (columns D = Price, E = Volume(1h)) with the
# latest scrape values, as produced by the "Updated cryptos list ... with
# GitHub Actions" job.
#
# Columns D/E are stored as plain text (e.g. "61.576.18", "0.999",
# "  -3.79%  ") rather than numbers, so any value that LOOKS like a plain
# number (single decimal point, e.g. "554.38") is written with a leading
# apostrophe to force Excel to keep it as text instead of silently
# reinterpreting it as a numeric value. The apostrophe flips the cell to a
# "quote prefix" style, so we immediately reset the cell style back to
# Normal afterwards to avoid leaving a stray formatting change behind -
# only the text content should differ from the original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Value
    )

    # Multi-dot "thousands" style values (e.g. "61.576.18") and values
    # containing other non-numeric characters (e.g. "  -3.79%  ") are
    # already safe - Excel can't coerce them to a number - so they can be
    # assigned directly and keep the default (style 0) formatting.
    if ($Value -match '^-?\d+(\.\d+)?$') {
        $ws.Range($Cell).Value = "'" + $Value
        $ws.Range($Cell).Style = "Normal"
    } else {
        $ws.Range($Cell).Value = $Value
    }
}

Set-TextValue "D2"  "61.576.18"
Set-TextValue "E2"  "  -3.79%  "

Set-TextValue "D3"  "2.476.76"
Set-TextValue "E3"  "  -6.40%  "

Set-TextValue "E4"  "  +0.03%  "

Set-TextValue "D5"  "554.38"
Set-TextValue "E5"  "  -4.86%  "

Set-TextValue "D6"  "147.09"
Set-TextValue "E6"  "  -5.84%  "

Set-TextValue "E7"  "  +0.03%  "

Set-TextValue "E8"  "  -3.51%  "

Set-TextValue "D9"  "2.473.34"
Set-TextValue "E9"  "  -6.41%  "

Set-TextValue "E10" "  -9.10%  "

Set-TextValue "D11" "5.47"
Set-TextValue "E11" "  -6.03%  "

Set-TextValue "E12" "  -1.44%  "

Set-TextValue "D13" "0.357"
Set-TextValue "E13" "  -6.88%  "

Set-TextValue "D14" "26.18"
Set-TextValue "E14" "  -8.21%  "

Set-TextValue "D15" "2.921.64"
Set-TextValue "E15" "  -6.40%  "

Set-TextValue "D16" "0.0000168"
Set-TextValue "E16" "  -9.59%  "

Set-TextValue "D17" "61.427.80"
Set-TextValue "E17" "  -3.82%  "

Set-TextValue "D18" "2.480.59"
Set-TextValue "E18" "  -6.32%  "

Set-TextValue "D19" "11.21"
Set-TextValue "E19" "  -8.11%  "

Set-TextValue "D20" "7.04"
Set-TextValue "E20" "  -9.08%  "

Set-TextValue "D21" "4.20"
Set-TextValue "E21" "  -7.52%  "

Set-TextValue "D22" "321.93"
Set-TextValue "E22" "  -7.04%  "

Set-TextValue "E23" "  -0.04%  "

Set-TextValue "E24" "  -0.80%  "

Set-TextValue "D25" "64.11"
Set-TextValue "E25" "  -5.79%  "

Set-TextValue "D26" "0.0000100"
Set-TextValue "E26" "  -11.30%  "

Set-TextValue "D27" "2.600.95"
Set-TextValue "E27" "  -6.13%  "

Set-TextValue "E28" "  -6.90%  "

Set-TextValue "D29" "546.22"
Set-TextValue "E29" "  -11.33%  "

Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  -0.02%  "

Set-TextValue "D31" "8.31"
Set-TextValue "E31" "  -10.52%  "

Set-TextValue "D32" "7.72"
Set-TextValue "E32" "  -5.58%  "

Set-TextValue "D33" "0.151"
Set-TextValue "E33" "  -6.53%  "

Set-TextValue "D34" "1.93"
Set-TextValue "E34" "  -7.09%  "

Set-TextValue "E35" "  -9.04%  "

Set-TextValue "D36" "5.92"
Set-TextValue "E36" "  -10.40%  "

Set-TextValue "E37" "  -10.70%  "

Set-TextValue "E38" "  +0.01%  "

Set-TextValue "D39" "0.383"
Set-TextValue "E39" "  -5.23%  "

Set-TextValue "D40" "18.58"
Set-TextValue "E40" "  -5.96%  "

Set-TextValue "D41" "145.23"
Set-TextValue "E41" "  -4.11%  "

Set-TextValue "D42" "1.74"
Set-TextValue "E42" "  -8.87%  "

Set-TextValue "D44" "40.42"
Set-TextValue "E44" "  -3.64%  "

Set-TextValue "D45" "2.36"
Set-TextValue "E45" "  -7.72%  "

Set-TextValue "D46" "147.95"
Set-TextValue "E46" "  -9.17%  "

Set-TextValue "D47" "3.65"
Set-TextValue "E47" "  -6.61%  "

Set-TextValue "D48" "21.73"
Set-TextValue "E48" "  -10.36%  "

Set-TextValue "E49" "  -8.41%  "

Set-TextValue "D50" "0.596"
Set-TextValue "E50" "  -6.28%  "

Set-TextValue "D51" "0.0944"
Set-TextValue "E51" "  -5.50%  "
